$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Column A (Mat/NC) - numbers, no new shared strings
$ws.Cells.Item(2, 1).Value = 21330051920003
$ws.Cells.Item(3, 1).Value = 21330051920017

# Column B (Paterno) - new shared strings BAEZ, MARTINEZ
$ws.Cells.Item(2, 2).Value = "BAEZ"
$ws.Cells.Item(3, 2).Value = "MARTINEZ"

# Column C (Materno) - new shared strings MARCELINO, XOTLANIHUA
$ws.Cells.Item(2, 3).Value = "MARCELINO"
$ws.Cells.Item(3, 3).Value = "XOTLANIHUA"

# Column D (Nombres) - new shared strings LUIS EDUARDO, YAIR
$ws.Cells.Item(2, 4).Value = "LUIS EDUARDO"
$ws.Cells.Item(3, 4).Value = "YAIR"

# Column E (Nombre_Largo) - existing shared string ÁLGEBRA
$ws.Cells.Item(2, 5).Value = "ÁLGEBRA"
$ws.Cells.Item(3, 5).Value = "ÁLGEBRA"

# Column F (Grupo) - existing shared string 1AV
$ws.Cells.Item(2, 6).Value = "1AV"
$ws.Cells.Item(3, 6).Value = "1AV"

# Column G (Reprobadas) - numbers
$ws.Cells.Item(2, 7).Value = 6
$ws.Cells.Item(3, 7).Value = 6
